$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1835.75
$ws.Range("I40").Value = 1868.3438
$ws.Range("J40").Value = 1575
$ws.Range("K40").Value = 1868.3438
$ws.Range("L40").Value = 1575
$ws.Range("M40").Value = -1693.3438
$ws.Range("N40").Value = -1925

$ws.Range("H43").Value = 1199
$ws.Range("I43").Value = 1199
$ws.Range("K43").Value = 1199
$ws.Range("M43").Value = -1130

$ws.Range("H80").Value = 1757.4
$ws.Range("I80").Value = 874.8
$ws.Range("K80").Value = 2624.4
$ws.Range("M80").Value = -1626.4

$ws.Range("H83").Value = 1757.4
$ws.Range("I83").Value = 874.8
$ws.Range("K83").Value = 7873.2
$ws.Range("M83").Value = -2881.2

$ws.Range("H96").Value = 17966.5
$ws.Range("I96").Value = 25574.75
$ws.Range("K96").Value = 76724.25
$ws.Range("M96").Value = -75351.25

$ws.Range("H132").Value = 2695.7273
$ws.Range("I132").Value = 2881.125
$ws.Range("J132").Value = 2201.3333
$ws.Range("K132").Value = 8643.375
$ws.Range("L132").Value = 6603.999899999999
$ws.Range("M132").Value = -6113.375
$ws.Range("N132").Value = -11663.9999

$ws.Range("H138").Value = 5321.846
$ws.Range("J138").Value = 5771.9443
$ws.Range("L138").Value = 17315.8329
$ws.Range("N138").Value = -27595.8329

$ws.Range("H141").Value = 2260.111
$ws.Range("I141").Value = 1140.3334
$ws.Range("K141").Value = 3421.0002
$ws.Range("M141").Value = 1758.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2752566.2
$ws.Range("I32").Value = 2595253.8
$ws.Range("K32").Value = 2595253.8
$ws.Range("M32").Value = -2594966.8

$ws.Range("H36").Value = 14342
$ws.Range("I36").Value = 4013
$ws.Range("J36").Value = 35000
$ws.Range("K36").Value = 4013
$ws.Range("L36").Value = 35000
$ws.Range("M36").Value = -3667
$ws.Range("N36").Value = -35692

$ws.Range("H61").Value = 3091.6667
$ws.Range("I61").Value = 3091.6667
$ws.Range("K61").Value = 3091.6667
$ws.Range("M61").Value = -2879.6667

$ws.Range("H132").Value = 1499
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 4497
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -9557

$ws.Range("H136").Value = 3091.6667
$ws.Range("I136").Value = 3091.6667
$ws.Range("K136").Value = 9275.000100000001
$ws.Range("M136").Value = -6725.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 167000
$ws.Range("J4").Value = 167000
$ws.Range("L4").Value = 167000
$ws.Range("N4").Value = -167224

$ws.Range("H31").Value = 2347.3438
$ws.Range("I31").Value = 1167.2858
$ws.Range("K31").Value = 1167.2858
$ws.Range("M31").Value = -872.2858000000001

$ws.Range("H34").Value = 2347.3438
$ws.Range("I34").Value = 1167.2858
$ws.Range("K34").Value = 1167.2858
$ws.Range("M34").Value = -965.2858000000001

$ws.Range("H122").Value = 921.4167
$ws.Range("I122").Value = 1000.8
$ws.Range("J122").Value = 524.5
$ws.Range("K122").Value = 3002.4
$ws.Range("L122").Value = 1573.5
$ws.Range("M122").Value = -552.3999999999996
$ws.Range("N122").Value = -6473.5

$ws.Range("H132").Value = 2720.3333
$ws.Range("I132").Value = 2720.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8160.999899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5630.999899999999
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 2757.1
$ws.Range("I134").Value = 2595.1667
$ws.Range("K134").Value = 7785.500100000001
$ws.Range("M134").Value = -5250.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

$ws.Range("H68").Value = 4999.8335
$ws.Range("J68").Value = 4999.8335
$ws.Range("L68").Value = 14999.5005
$ws.Range("N68").Value = -16621.5005

$ws.Range("H71").Value = 4999.8335
$ws.Range("J71").Value = 4999.8335
$ws.Range("L71").Value = 44998.5015
$ws.Range("N71").Value = -53110.5015

$ws.Range("H131").Value = 835051
$ws.Range("I131").Value = 1225.8572
$ws.Range("K131").Value = 3677.5716
$ws.Range("M131").Value = 1362.4284

$ws.Range("H132").Value = 1958.6
$ws.Range("I132").Value = 1632.6666
$ws.Range("J132").Value = 2447.5
$ws.Range("K132").Value = 14693.9994
$ws.Range("L132").Value = 22027.5
$ws.Range("M132").Value = -12163.9994
$ws.Range("N132").Value = -27087.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 71357.14
$ws.Range("I57").Value = 22000
$ws.Range("J57").Value = 79583.336
$ws.Range("K57").Value = 22000
$ws.Range("L57").Value = 79583.336
$ws.Range("M57").Value = -21180
$ws.Range("N57").Value = -81223.336

$ws.Range("H97").Value = 750
$ws.Range("I97").Value = 750
$ws.Range("K97").Value = 750
$ws.Range("M97").Value = -254

$ws.Range("H126").Value = 3051.6
$ws.Range("I126").Value = 1275
$ws.Range("J126").Value = 4236
$ws.Range("K126").Value = 3825
$ws.Range("L126").Value = 12708
$ws.Range("M126").Value = -1355
$ws.Range("N126").Value = -17648

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1680.5834
$ws.Range("I46").Value = 1533.375
$ws.Range("J46").Value = 1975
$ws.Range("K46").Value = 1533.375
$ws.Range("L46").Value = 1975
$ws.Range("M46").Value = -1345.375
$ws.Range("N46").Value = -2351

$ws.Range("H58").Value = 4700
$ws.Range("I58").Value = 4700
$ws.Range("K58").Value = 4700
$ws.Range("M58").Value = -4440

$ws.Range("H61").Value = 1888
$ws.Range("I61").Value = 1623.5
$ws.Range("J61").Value = 2064.3333
$ws.Range("K61").Value = 1623.5
$ws.Range("L61").Value = 2064.3333
$ws.Range("M61").Value = -1421.5
$ws.Range("N61").Value = -2468.3333

$ws.Range("H68").Value = 1851.2727
$ws.Range("I68").Value = 1608.25
$ws.Range("J68").Value = 2499.3333
$ws.Range("K68").Value = 1608.25
$ws.Range("L68").Value = 2499.3333
$ws.Range("M68").Value = -859.25
$ws.Range("N68").Value = -3997.3333

$ws.Range("H71").Value = 1851.2727
$ws.Range("I71").Value = 1608.25
$ws.Range("J71").Value = 2499.3333
$ws.Range("K71").Value = 8041.25
$ws.Range("L71").Value = 12496.6665
$ws.Range("M71").Value = -4297.25
$ws.Range("N71").Value = -19984.6665

$ws.Range("H93").Value = 1094.6666
$ws.Range("I93").Value = 1193.6
$ws.Range("J93").Value = 600
$ws.Range("K93").Value = 1193.6
$ws.Range("L93").Value = 600
$ws.Range("M93").Value = 54.40000000000009
$ws.Range("N93").Value = -3096

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H113").Value = 1888
$ws.Range("I113").Value = 1623.5
$ws.Range("J113").Value = 2064.3333
$ws.Range("K113").Value = 1623.5
$ws.Range("L113").Value = 2064.3333
$ws.Range("M113").Value = 546.5
$ws.Range("N113").Value = -6404.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 455090.53
$ws.Range("I5").Value = 599.6
$ws.Range("K5").Value = 599.6
$ws.Range("M5").Value = -487.6

$ws.Range("H59").Value = 38000
$ws.Range("J59").Value = 38000
$ws.Range("L59").Value = 38000
$ws.Range("N59").Value = -39476

$ws.Range("H132").Value = 509.66666
$ws.Range("I132").Value = 610.6
$ws.Range("J132").Value = 5
$ws.Range("K132").Value = 1831.8
$ws.Range("L132").Value = 15
$ws.Range("M132").Value = 698.1999999999998
$ws.Range("N132").Value = -5075

$ws.Range("H136").Value = 2822.6924
$ws.Range("I136").Value = 2538.6956
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7616.0868
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5066.0868
$ws.Range("N136").Value = -20100
